$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 39, pushing existing rows 39:54 down to 40:55.
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with the new weekly record.
$ws.Cells.Item(39, 1).Value = 1
$ws.Cells.Item(39, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(39, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(39, 4).Value = 44489
$ws.Cells.Item(39, 5).Value = 15
$ws.Cells.Item(39, 6).Value = 100112040
$ws.Cells.Item(39, 7).Value = "Cilantro"
$ws.Cells.Item(39, 8).Value = "Sin especificar"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 250
$ws.Cells.Item(39, 11).Value = 900
$ws.Cells.Item(39, 12).Value = 1000
$ws.Cells.Item(39, 13).Value = 950
$ws.Cells.Item(39, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(39, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(39, 16).Value = 475
$ws.Cells.Item(39, 17).Value = 2
$ws.Cells.Item(39, 18).Value = "Hortaliza"
